$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $old"
    }
}

# 1) Merge the split "docker==18.09.5 ... installed):" run fragments into one run.
Replace-Text "need to have docker==18.09.5 and docker-compose==1.23.2 installed):" "need to have docker==18.09.5 and docker-compose==1.23.2 installed):"

# 2) "Clone Flapweb's repo and start docker container" -> capitalize Docker
Replace-Text "Clone Flapweb’s repo and start docker container" "Clone Flapweb’s repo and start Docker container"

# 3) Merge "Install requirements: " + "docker-compose build"
Replace-Text "Install requirements: docker-compose build" "Install requirements: docker-compose build"

# 4) Merge "(optional) Start app’s container: " + "docker-compose up"
Replace-Text "(optional) Start app’s container: docker-compose up" "(optional) Start app’s container: docker-compose up"

# 5) Merge "Enter to the docker container env: " + "docker exec -it app_flap bash"
Replace-Text "Enter to the docker container env: docker exec -it app_flap bash" "Enter to the docker container env: docker exec -it app_flap bash"

# 6) Merge "Start postgresql app: " + "psql -U guillermo -h pg_flap"
Replace-Text "Start postgresql app: psql -U guillermo -h pg_flap" "Start postgresql app: psql -U guillermo -h pg_flap"

# 7) Merge "Create database: " + "CREATE DATABASE registry;"
Replace-Text "Create database: CREATE DATABASE registry;" "Create database: CREATE DATABASE registry;"

# 8) Merge "Make migrations" + ":"
Replace-Text "Make migrations:" "Make migrations:"

# 9) Merge "You can skip the last step by " + "quitting" + " psql."
Replace-Text "You can skip the last step by quitting psql." "You can skip the last step by quitting psql."

# 10) Merge "Migrate Django framework database tables: " + "python manage.py migrate"
Replace-Text "Migrate Django framework database tables: python manage.py migrate" "Migrate Django framework database tables: python manage.py migrate"

# 11) Merge "Create the migrations (generate the " + "SQL " + "commands): "
Replace-Text "Create the migrations (generate the SQL commands): " "Create the migrations (generate the SQL commands): "

# 12) Merge "Run the migrations (execute the SQL commands): " + "python manage.py migrate"
Replace-Text "Run the migrations (execute the SQL commands): python manage.py migrate" "Run the migrations (execute the SQL commands): python manage.py migrate"

# 13) Merge "Start app’s container: " + "docker-compose up" (the one near bookmarkEnd, without "(optional)")
Replace-Text "Start app’s container: docker-compose up" "Start app’s container: docker-compose up"

# 14) Add the new character style "ListLabel 3" / styleId ListLabel3
$newStyle = $d.Styles.Add("ListLabel 3", 2)
$newStyle.QuickStyle = $true

Write-Host "done"
